# "Final Build, with & w/o JSE, added Screenshots, deleted comments"
#
# The "Doctors" sheet's sample-data block (rows 2-6, columns A:D) is
# refreshed with a new batch of doctors (Kolkata-based ENT specialists),
# replacing the previous ("Homoeopath", Pune-based) sample rows. The
# "Surgeries List" column (E) and everything on the "Demo" sheet is left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Doctors")

# Row 2
$ws.Cells.Item(2, 1).Value = "Dr. Pranay Kumar Agarwal"
$ws.Cells.Item(2, 2).Value = "Ear-Nose-Throat (ENT) Specialist"
$ws.Cells.Item(2, 3).Value = "13 years experience overall"
$ws.Cells.Item(2, 4).Value = "Park Circus,Kolkata"

# Row 3
$ws.Cells.Item(3, 1).Value = "Dr. Sudipta Chandra"
$ws.Cells.Item(3, 2).Value = "Ear-Nose-Throat (ENT) Specialist"
$ws.Cells.Item(3, 3).Value = "24 years experience overall"
$ws.Cells.Item(3, 4).Value = "Minto Park,Kolkata"

# Row 4
$ws.Cells.Item(4, 1).Value = "Dr. Sayan Ganguly"
$ws.Cells.Item(4, 2).Value = "Ear-Nose-Throat (ENT) Specialist"
$ws.Cells.Item(4, 3).Value = "31 years experience overall"
$ws.Cells.Item(4, 4).Value = "Minto Park,Kolkata"

# Row 5
$ws.Cells.Item(5, 1).Value = "Dr. Sunil Jalan"
$ws.Cells.Item(5, 2).Value = "Ear-Nose-Throat (ENT) Specialist"
$ws.Cells.Item(5, 3).Value = "20 years experience overall"
$ws.Cells.Item(5, 4).Value = "Sarat Bose Road,Kolkata"

# Row 6
$ws.Cells.Item(6, 1).Value = "Dr. Nitin Mittal"
$ws.Cells.Item(6, 2).Value = "Ear-Nose-Throat (ENT) Specialist"
$ws.Cells.Item(6, 3).Value = "18 years experience overall"
$ws.Cells.Item(6, 4).Value = "Lake Town,Kolkata"
